$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 77: correct the date and record the payment received for the
# existing "Eko 4B" listrik entry (Dibayar column J + new Diterima/H value)
$ws.Range("A77").Value = 43211
$ws.Range("H77").NumberFormat = $ws.Range("F77").NumberFormat
$ws.Range("H77").Value = 51500
$ws.Range("J77").Value = 51500

# Row 78: new incoming payment - Bagus Prasojo (Listrik)
$ws.Range("A78").Value = 43211
$ws.Range("B78").Value = "BL171111CBC2ELC"
$ws.Range("C78").Value = "Bagus Prasojo"
$ws.Range("D78").Value = "520530520059"
$ws.Range("D78").NumberFormat = "@"
$ws.Range("E78").Value = "Bagus Prasojo"
$ws.Range("E78").NumberFormat = "@"
$ws.Range("F78").Value = 51500
$ws.Range("G78").Value = "Listrik"
$ws.Range("H78").NumberFormat = $ws.Range("F77").NumberFormat
$ws.Range("H78").Value = 51500
$ws.Range("J78").Value = 51500

# Row 79: new incoming payment - Srini (Pulsa), saldo masuk dari Abudullah Tuyar
$ws.Range("C79").Value = "Srini"
$ws.Range("D79").Value = "08522904326"
$ws.Range("D79").NumberFormat = "@"
$ws.Range("E79").Value = "Srini"
$ws.Range("E79").NumberFormat = "@"
$ws.Range("F79").Value = 11000
$ws.Range("G79").Value = "Pulsa"
$ws.Range("H79").Value = 11000
$ws.Range("J79").Value = 11000

$ws.Range("J81").Select() | Out-Null
